# ------------------------------------------------------------------
# t2 - copia.xlsx : add the "Fibonacci" variants of the three tables
# (2^10 / 2^12 / 2^14), keep the original three sheets but rename
# them to "Heap ..." and re-colour all the tabs.
#
# NOTE: worksheet object references handed out by this host appear to
# be index-bound "live views" rather than stable handles - once a
# Move/Copy reshuffles tab positions, an old variable can silently
# start pointing at a different sheet. To stay safe we always re-fetch
# sheets by their (stable) Name right before using them.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# The workbook currently has 3 sheets: "2^10", "2^12", "2^14".
# "2^14" is the one that already carries the extra styled cell (B2)
# that the new "Fibonacci" sheets are based on, so we duplicate it
# three times, inserting each copy right before "2^10" (which keeps
# the original three sheets together and in their original order).

$template = $wb.Worksheets.Item("2^14")
$template.Copy($wb.Worksheets.Item("2^10"))
$wb.Worksheets.Item("2^14 (2)").Name = "Fibonacci 2^10"

$wb.Worksheets.Item("Fibonacci 2^10").Copy($wb.Worksheets.Item("2^10"))
$wb.Worksheets.Item("Fibonacci 2^10 (2)").Name = "Fibonacci 2^12"

$wb.Worksheets.Item("Fibonacci 2^12").Copy($wb.Worksheets.Item("2^10"))
$wb.Worksheets.Item("Fibonacci 2^12 (2)").Name = "Fibonacci 2^14"

# --- rename the original sheets -----------------------------------------
$wb.Worksheets.Item("2^10").Name = "Heap 2^10"
$wb.Worksheets.Item("2^12").Name = "Heap 2^12"
$wb.Worksheets.Item("2^14").Name = "Heap 2^14"

# --- tab colours ----------------------------------------------------------
# "Fibonacci" tabs -> Orange, Accent 2 (theme 5, no tint) = RGB(237,125,49)
$wb.Worksheets.Item("Fibonacci 2^10").Tab.Color = 3243501
$wb.Worksheets.Item("Fibonacci 2^12").Tab.Color = 3243501
$wb.Worksheets.Item("Fibonacci 2^14").Tab.Color = 3243501

# "Heap" tabs -> Green, Accent 6, Lighter 40% (theme 9, tint ~0.4) = RGB(169,209,142)
$wb.Worksheets.Item("Heap 2^10").Tab.Color = 9359785
$wb.Worksheets.Item("Heap 2^12").Tab.Color = 9359785
$wb.Worksheets.Item("Heap 2^14").Tab.Color = 9359785

# --- Fibonacci 2^14 gets an extra formatted (but empty) cell at J31 -------
$wb.Worksheets.Item("Fibonacci 2^14").Range("J31").Font.Underline = $true

# --- selections on every sheet (also drives which tab ends up "active") --
$wb.Worksheets.Item("Heap 2^10").Activate()
$wb.Worksheets.Item("Heap 2^10").Range("B1").Select()

$wb.Worksheets.Item("Heap 2^12").Activate()
$wb.Worksheets.Item("Heap 2^12").Range("E35").Select()

$wb.Worksheets.Item("Heap 2^14").Activate()
$wb.Worksheets.Item("Heap 2^14").Range("H37").Select()

$wb.Worksheets.Item("Fibonacci 2^10").Activate()
$wb.Worksheets.Item("Fibonacci 2^10").Range("H25").Select()

$wb.Worksheets.Item("Fibonacci 2^12").Activate()
$wb.Worksheets.Item("Fibonacci 2^12").Range("E34").Select()

# "Fibonacci 2^14" is left as the active / selected sheet, matching the
# workbook's activeTab pointing at it.
$wb.Worksheets.Item("Fibonacci 2^14").Activate()
$wb.Worksheets.Item("Fibonacci 2^14").Range("J31").Select()
